$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.993294715881348
$ws.Range("B1").Value = 0.2977640330791473
$ws.Range("C1").Value = 0.2140259593725204
$ws.Range("D1").Value = 0.2033171355724335
$ws.Range("E1").Value = 0.2152483463287354
